$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = "'29.857.10"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value2 = "  -0.06%  "
$ws.Cells.Item(3, 4).Value2 = "'1.887.43"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value2 = "  -0.39%  "
$ws.Cells.Item(4, 4).Value2 = "'1.000"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value2 = "  -0.01%  "
$ws.Cells.Item(5, 4).Value2 = "'0.7469"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value2 = "  -4.69%  "
$ws.Cells.Item(6, 4).Value2 = "'242.60"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value2 = "  -0.42%  "
$ws.Cells.Item(7, 5).Value2 = "  +0.00%  "
$ws.Cells.Item(8, 4).Value2 = "'0.3113"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value2 = "  -0.87%  "
$ws.Cells.Item(9, 5).Value2 = "  -1.80%  "
$ws.Cells.Item(10, 4).Value2 = "'0.07122"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value2 = "  -2.06%  "
$ws.Cells.Item(11, 4).Value2 = "'0.08488"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value2 = "  +4.77%  "
$ws.Cells.Item(12, 4).Value2 = "'0.7603"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value2 = "  -1.99%  "
$ws.Cells.Item(13, 4).Value2 = "'1.912.59"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value2 = "  -0.82%  "
$ws.Cells.Item(14, 4).Value2 = "'5.357"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value2 = "  -2.33%  "
$ws.Cells.Item(15, 4).Value2 = "'93.38"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value2 = "  -0.89%  "
$ws.Cells.Item(16, 4).Value2 = "'6.156"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value2 = "  -0.86%  "
$ws.Cells.Item(17, 4).Value2 = "'29.897.17"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value2 = "  +0.24%  "
$ws.Cells.Item(18, 4).Value2 = "'13.70"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value2 = "  -1.81%  "
$ws.Cells.Item(19, 4).Value2 = "'243.32"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value2 = "  -1.30%  "
$ws.Cells.Item(20, 4).Value2 = "'0.000007794"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value2 = "  -0.37%  "
$ws.Cells.Item(21, 4).Value2 = "'2.159.68"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value2 = "  +3.26%  "
$ws.Cells.Item(22, 5).Value2 = "  +0.03%  "
$ws.Cells.Item(23, 4).Value2 = "'7.983"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value2 = "  -1.62%  "
$ws.Cells.Item(24, 4).Value2 = "'1.001"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(25, 4).Value2 = "'0.1589"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value2 = "  -0.10%  "
$ws.Cells.Item(26, 4).Value2 = "'9.380"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value2 = "  -0.85%  "
$ws.Cells.Item(27, 4).Value2 = "'162.30"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value2 = "  -1.10%  "
$ws.Cells.Item(28, 4).Value2 = "'18.76"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value2 = "  -0.02%  "
$ws.Cells.Item(29, 4).Value2 = "'2.027"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value2 = "  +0.15%  "
$ws.Cells.Item(30, 4).Value2 = "'1.502"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value2 = "  +4.22%  "
$ws.Cells.Item(31, 4).Value2 = "'1.530"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value2 = "  -0.83%  "
$ws.Cells.Item(32, 4).Value2 = "'4.474"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value2 = "  +0.05%  "
$ws.Cells.Item(33, 4).Value2 = "'4.102"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value2 = "  +0.69%  "
$ws.Cells.Item(34, 4).Value2 = "'0.05390"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value2 = "  -3.19%  "
$ws.Cells.Item(35, 4).Value2 = "'1.239"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value2 = "  -0.40%  "
$ws.Cells.Item(36, 4).Value2 = "'0.7440"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value2 = "  -1.32%  "
$ws.Cells.Item(37, 4).Value2 = "'1.001"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value2 = "  -0.10%  "
$ws.Cells.Item(38, 4).Value2 = "'2.710"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value2 = "  +1.22%  "
$ws.Cells.Item(40, 5).Value2 = "  -0.96%  "
$ws.Cells.Item(41, 4).Value2 = "'0.4458"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value2 = "  -0.06%  "
$ws.Cells.Item(42, 4).Value2 = "'6.059"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value2 = "  +1.75%  "
$ws.Cells.Item(43, 4).Value2 = "'72.71"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value2 = "  -2.11%  "
$ws.Cells.Item(44, 4).Value2 = "'1.088.67"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value2 = "  -4.64%  "
$ws.Cells.Item(45, 4).Value2 = "'0.8590"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value2 = "  +0.73%  "
$ws.Cells.Item(46, 4).Value2 = "'1.000"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value2 = "  +0.03%  "
$ws.Cells.Item(47, 4).Value2 = "'102.34"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value2 = "  +0.31%  "
$ws.Cells.Item(48, 4).Value2 = "'7.673"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value2 = "  +1.70%  "
$ws.Cells.Item(49, 4).Value2 = "'1.862"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value2 = "  -1.69%  "
$ws.Cells.Item(50, 4).Value2 = "'3.069"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value2 = "  -2.41%  "
$ws.Cells.Item(51, 4).Value2 = "'2.054.10"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value2 = "  +1.22%  "
